{"js": "// Replace the arithmetic expressions in the practice-sheet table with a\n// new set of problems, cell-for-cell, in row-major (document) order.\n// Formatting (fonts, size, alignment, etc.) is untouched because we only\n// rewrite the cell text via the Table.values setter.\n\nconst newProblems = [\n  \"71-53=\", \"85-75=\", \"19-9=\", \"37-17=\", \"99-37=\",\n  \"70-33=\", \"45-28=\", \"78-43=\", \"56-29=\", \"33+54=\",\n  \"51-2=\", \"1+32=\", \"99-93=\", \"1+40=\", \"29+62=\",\n  \"47-30=\", \"6+43=\", \"99-68=\", \"43-2=\", \"91-16=\",\n  \"93-2=\", \"23+6=\", \"9-0=\", \"89-1=\", \"98-45=\",\n  \"59-23=\", \"17+31=\", \"97-31=\", \"61+10=\", \"4+64=\",\n  \"76-55=\", \"89-16=\", \"62-42=\", \"43-20=\", \"51+18=\",\n  \"42-1=\", \"91-38=\", \"6+70=\", \"71-35=\", \"66-34=\",\n  \"8+1=\", \"61-38=\", \"47-16=\", \"15+36=\", \"66-43=\",\n  \"5-1=\", \"86+13=\", \"73-6=\", \"13+55=\", \"90-7=\",\n  \"38+0=\", \"14+44=\", \"77-25=\", \"84+6=\", \"40-14=\",\n  \"48+34=\", \"25+26=\", \"71+8=\", \"6+9=\", \"15+43=\",\n  \"1+0=\", \"18-16=\", \"50-18=\", \"24+25=\", \"72-45=\",\n  \"4+30=\", \"85-71=\", \"75-75=\", \"64-63=\", \"57-20=\",\n  \"56+26=\", \"25-10=\", \"21+21=\", \"63+10=\", \"55-26=\",\n  \"54-32=\", \"23+31=\", \"33+55=\", \"93-52=\", \"88-64=\",\n  \"95-36=\", \"48+34=\", \"93-71=\", \"19+9=\", \"39-36=\",\n  \"40-18=\", \"16-5=\", \"35-32=\", \"20+13=\", \"61+28=\",\n  \"17+72=\", \"45+9=\", \"95-28=\", \"55-46=\", \"36+1=\",\n  \"82-67=\", \"12+42=\", \"46+37=\", \"78-60=\", \"83-59=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst colCount = table.values[0].length;\nconst newValues = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  const row = [];\n  for (let c = 0; c < colCount; c++) {\n    row.push(newProblems[r * colCount + c]);\n  }\n  newValues.push(row);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the arithmetic expressions in the practice-sheet table with a\n# new set of problems, cell-for-cell, in row-major (document) order.\n# Only the cell Range.Text is rewritten, so run formatting (fonts, size,\n# alignment, etc.) is left untouched.\n\n$newProblems = @(\n  @(\"71-53=\", \"85-75=\", \"19-9=\", \"37-17=\", \"99-37=\"),\n  @(\"70-33=\", \"45-28=\", \"78-43=\", \"56-29=\", \"33+54=\"),\n  @(\"51-2=\", \"1+32=\", \"99-93=\", \"1+40=\", \"29+62=\"),\n  @(\"47-30=\", \"6+43=\", \"99-68=\", \"43-2=\", \"91-16=\"),\n  @(\"93-2=\", \"23+6=\", \"9-0=\", \"89-1=\", \"98-45=\"),\n  @(\"59-23=\", \"17+31=\", \"97-31=\", \"61+10=\", \"4+64=\"),\n  @(\"76-55=\", \"89-16=\", \"62-42=\", \"43-20=\", \"51+18=\"),\n  @(\"42-1=\", \"91-38=\", \"6+70=\", \"71-35=\", \"66-34=\"),\n  @(\"8+1=\", \"61-38=\", \"47-16=\", \"15+36=\", \"66-43=\"),\n  @(\"5-1=\", \"86+13=\", \"73-6=\", \"13+55=\", \"90-7=\"),\n  @(\"38+0=\", \"14+44=\", \"77-25=\", \"84+6=\", \"40-14=\"),\n  @(\"48+34=\", \"25+26=\", \"71+8=\", \"6+9=\", \"15+43=\"),\n  @(\"1+0=\", \"18-16=\", \"50-18=\", \"24+25=\", \"72-45=\"),\n  @(\"4+30=\", \"85-71=\", \"75-75=\", \"64-63=\", \"57-20=\"),\n  @(\"56+26=\", \"25-10=\", \"21+21=\", \"63+10=\", \"55-26=\"),\n  @(\"54-32=\", \"23+31=\", \"33+55=\", \"93-52=\", \"88-64=\"),\n  @(\"95-36=\", \"48+34=\", \"93-71=\", \"19+9=\", \"39-36=\"),\n  @(\"40-18=\", \"16-5=\", \"35-32=\", \"20+13=\", \"61+28=\"),\n  @(\"17+72=\", \"45+9=\", \"95-28=\", \"55-46=\", \"36+1=\"),\n  @(\"82-67=\", \"12+42=\", \"46+37=\", \"78-60=\", \"83-59=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $t.Cell($r, $c).Range.Text = $newProblems[$r - 1][$c - 1]\n  }\n}\n"}
